$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, matching the style/formatting of the other header
# cells (B1:G1) by copying G1's format onto it.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New value in H2 (plain/default formatting, like the other data cells in row 2)
$ws.Range("H2").Value = 1
